# Insert three new bullet paragraphs right after the
# "GIS & Geospatial Analysis Consulting" paragraph under the
# PARTNER - Siege Analytics entry, before the existing "• Lead ..." bullet.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*GIS & Geospatial Analysis Consulting*") {
        $r = $p.Range
        $ip = $d.Range($r.End, $r.End)
        $ip.InsertAfter("• Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels`r")
        $ip.InsertAfter("• Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide`r")
        $ip.InsertAfter("• Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis`r")
        break
    }
}
